$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F10").Value = "Comments back!"
$ws.Range("E10").Value = "Needs rewrite"
$ws.Range("E8").Value = "Needs polish"
$ws.Range("E9").Value = "Needs polish"

$ws.Range("E11").Select()
